# Generate Report for Handoff
# Updates the "b.md" row across the Overview/zh-cn/de-de sheets to reflect
# that the file has been re-handed-off (new xlf files generated), and
# records the "version not latest" error detail on the handback row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$status        = "Ready for handoff"
$errorDetail   = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/e041d4e93ecc40c76a16a49b2e1063742e79ad44/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/54bb97d8806095e7337de9c70b5ea15780d1e898/e2e/b.md."

# ---- Overview sheet: row 3 corresponds to b.md ----
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = "2016-07-26 07:43:49"

# ---- zh-cn sheet: row 3 corresponds to b.md ----
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("F3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("G3").Value = "2016-07-26 07:43:38"
$wsZhCn.Range("O3").Value = $errorDetail
$wsZhCn.Range("O1:O3").ColumnWidth = 39.1667

# ---- de-de sheet: row 3 corresponds to b.md ----
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("F3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("G3").Value = "2016-07-26 07:43:49"
$wsDeDe.Range("O3").Value = $errorDetail
$wsDeDe.Range("O1:O3").ColumnWidth = 39.1667
